$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = 7630
$ws.Range("C2").Value = 4

# Delete rows 3 and 4 entirely (shrinks used range / dimension to A1:C2)
$ws.Range("A3:C4").Delete()
